$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.330500000000003
$ws.Range("B21").Value = 9.410600000000006
$ws.Range("B23").Value = 9.022999999999994
$ws.Range("B25").Value = 5.914399999999997
$ws.Range("B53").Value = 6.039800000000001
$ws.Range("B57").Value = 4.792099999999994
$ws.Range("B59").Value = 4.976399999999999
$ws.Range("B69").Value = 5.403499999999993
$ws.Range("B79").Value = 9.677900000000005
$ws.Range("B83").Value = 5.148799999999997
$ws.Range("B93").Value = 5.478199999999998
